# Applies the "Mas mediciones de temperatura" edit:
#  - Replaces the raw temperature/timestamp log on "Datos crudos" with a new
#    33-row run (rows 2-34), removing the old rows 35-36.
#  - Updates the manual pointer H2 (12 -> 16).
#  - Fixes up the selection/active-tab view state to match the authored file.
# Everything else (H3 count, I2 pointer delta, and all of "Datos validos")
# is formula-driven and recalculates automatically.

$wb = $excel.ActiveWorkbook
$wsRaw = $wb.Worksheets.Item("Datos crudos")
$wsValid = $wb.Worksheets.Item("Datos válidos")

$data = @(
    @{Row=2; C="2023-12-12 03:19:55"; E=26.984265734265701},
    @{Row=3; C="2023-12-12 03:20:56"; E=26.4597902097901},
    @{Row=4; C="2023-12-12 03:21:57"; E=25.935314685314601},
    @{Row=5; C="2023-12-12 03:22:58"; E=25.541958041958001},
    @{Row=6; C="2023-12-12 03:23:59"; E=25.148601398601301},
    @{Row=7; C="2023-12-12 03:25:00"; E=25.017482517482499},
    @{Row=8; C="2023-12-12 03:26:02"; E=24.624125874125799},
    @{Row=9; C="2023-12-12 03:27:03"; E=24.7552447552447},
    @{Row=10; C="2023-12-12 03:28:04"; E=24.361888111888099},
    @{Row=11; C="2023-12-12 03:29:05"; E=24.361888111888099},
    @{Row=12; C="2023-12-12 03:30:06"; E=24.230769230769202},
    @{Row=13; C="2023-12-12 03:31:08"; E=24.0996503496503},
    @{Row=14; C="2023-12-12 03:32:09"; E=24.0996503496503},
    @{Row=15; C="2023-12-12 03:33:10"; E=23.968531468531399},
    @{Row=16; C="2023-12-12 03:34:11"; E=23.837412587412501},
    @{Row=17; C="2023-12-12 03:35:12"; E=23.968531468531399},
    @{Row=18; C="2023-12-12 03:36:14"; E=23.968531468531399},
    @{Row=19; C="2023-12-12 03:37:15"; E=23.968531468531399},
    @{Row=20; C="2023-12-12 03:38:16"; E=23.575174825174798},
    @{Row=21; C="2023-12-12 03:39:17"; E=23.837412587412501},
    @{Row=22; C="2023-12-12 03:40:18"; E=23.837412587412501},
    @{Row=23; C="2023-12-12 03:41:19"; E=23.837412587412501},
    @{Row=24; C="2023-12-12 03:42:21"; E=23.7062937062937},
    @{Row=25; C="2023-12-12 03:43:22"; E=23.968531468531399},
    @{Row=26; C="2023-12-12 03:44:23"; E=23.837412587412501},
    @{Row=27; C="2023-12-12 03:45:24"; E=23.837412587412501},
    @{Row=28; C="2023-12-12 03:46:25"; E=23.444055944055901},
    @{Row=29; C="2023-12-12 03:47:27"; E=23.575174825174798},
    @{Row=30; C="2023-12-12 03:48:28"; E=23.575174825174798},
    @{Row=31; C="2023-12-12 03:49:29"; E=23.575174825174798},
    @{Row=32; C="2023-12-12 03:50:30"; E=23.7062937062937},
    @{Row=33; C="2023-12-12 03:51:31"; E=23.968531468531399},
    @{Row=34; C="2023-12-12 03:52:33"; E=23.837412587412501}
)

foreach ($row in $data) {
    $wsRaw.Cells.Item($row.Row, 3).Value = $row.C
    $wsRaw.Cells.Item($row.Row, 5).Value = $row.E
}

# Old data ran through row 36; the new run only has 33 rows (2-34), so the
# trailing two rows of the previous run are removed entirely.
$wsRaw.Rows.Item("35:36").Delete() | Out-Null

# Manual pointer used by I2 ( =INDEX(C:C,H2)-C2 )
$wsRaw.Range("H2").Value = 16

# View state: "Datos crudos" selection moves to G32 and loses the tab focus,
# "Datos válidos" becomes the active tab with selection M19.
$wsRaw.Range("G32").Select() | Out-Null
$wsValid.Range("M19").Select() | Out-Null
$wsValid.Activate() | Out-Null

$wb.Application.Calculate() | Out-Null
